# VTiger.xlsx update: add "Expected Title" column (B) with expected page
# titles to the Organization and Contacts sheets.

$wb = $excel.ActiveWorkbook

# --- Organization sheet (sheet1) ---------------------------------------
$ws1 = $wb.Worksheets.Item("Organization")

# Populate column B (order matters: it controls shared-string index order)
$ws1.Range("B2").Value = "Administrator - Organizations - vtiger CRM 5 - Commercial Open Source CRM"
$ws1.Range("B1").Value = "Expected Title"

# Style B2 with the Courier New font used to show the raw HTML <title>
$fontB2 = $ws1.Range("B2").Font
$fontB2.Name = "Courier New"
$fontB2.Size = 9
$fontB2.Color = 1972768

# Column B width ~= 75 characters
$ws1.Columns.Item(2).ColumnWidth = 74.16666666666667

# Selection + print orientation as left by the author after editing
$ws1.Range("B1").Select() | Out-Null
$ws1.PageSetup.Orientation = 1

# --- Contacts sheet (sheet2) --------------------------------------------
$ws2 = $wb.Worksheets.Item("Contacts")

$ws2.Range("B1").Value = "Expected Title"
$ws2.Range("B2").Value = " Administrator - Contacts - vtiger CRM 5 - Commercial Open Source CRM"

# Column B width ~= 64.86 characters
$ws2.Columns.Item(2).ColumnWidth = 64.02213541666667

# Author left the whole column selected on this (active) sheet
$ws2.Range("B:B").Select() | Out-Null
